$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1) "Role" section: replace the "IF <Role>" trigger-phrase cells with
#        the actual S2 screener-code conditions --------------------------
$d.Content.Find.Execute("IF Physician", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "S2=1 OR S2=2 OR S2=3 OR S2=4 OR S2=5", 2) | Out-Null

$d.Content.Find.Execute("IF NP/PA", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "S2=6 OR S2=7", 2) | Out-Null

# --- 2) "Segments" section: drop the "from list" suffix and clear the
#        italic placeholder styling on the four Segment rows -------------
function Update-SegmentCell($rowIndex, $oldText, $newText) {
    $cell = $t.Cell($rowIndex, 2)
    $rng = $cell.Range
    $find = $rng.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Replacement.Font.Italic = 0
    $find.Execute($oldText, $true, $false, $false, $false, $false, `
                  $true, 1, $true, $newText, 2) | Out-Null
}

Update-SegmentCell 19 "Segment A from list" "Segment A"
Update-SegmentCell 20 "Segment B from list" "Segment B"
Update-SegmentCell 21 "Segment C from list" "Segment C"
Update-SegmentCell 22 "Segment D from list" "Segment D"
